# "1st changes of mifos to finflux"
# Insert a new (blank) column before column N ("Late") on the
# "Repayment Schedule" sheet, which shifts the existing N/O/P columns
# one position to the right (N->O, O->P, P->Q), and switches the active
# sheet/selection from "Transactions" to "Repayment Schedule".

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column at position N (column 14), shifting N:P -> O:Q
$wsSchedule.Columns.Item(14).Insert()

# Make "Repayment Schedule" the active sheet/tab and set its selection
$wsSchedule.Activate()
$wsSchedule.Range("S5").Select()
